$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values
$ws.Range("B2").Value = "None-MacBookPro"
$ws.Range("D2").Value = "Apple Inc"
$ws.Range("K2").Value = "Stocked"
$ws.Range("M2").Value = ""

# Add new row 3
$ws.Range("A3").Value = "Computer"
$ws.Range("B3").Value = "Seba Salgado-Latitude"
$ws.Range("C3").Value = "Rack A"
$ws.Range("D3").Value = "Dell inc."
$ws.Range("E3").Value = "Latitude"
$ws.Range("F3").Value = "CS08BY3"
$ws.Range("G3").Value = ""
$ws.Range("H3").Value = "Check"
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = "Stocked"
$ws.Range("L3").Value = ""
$ws.Range("M3").Value = "Seba Salgado"
$ws.Range("N3").Value = "Laptop"
